$d = $word.ActiveDocument

$newValues = @(
    "38+18=",
    "90+6=",
    "58-44=",
    "3+12=",
    "31+34=",
    "30+69=",
    "71-66=",
    "69-49=",
    "54-46=",
    "89-50=",
    "26+26=",
    "23+65=",
    "7+41=",
    "57-52=",
    "12+39=",
    "1+73=",
    "57-36=",
    "49-0=",
    "62+1=",
    "76-69=",
    "21+38=",
    "92-27=",
    "9-3=",
    "9+85=",
    "65-50=",
    "6+30=",
    "19+50=",
    "35+53=",
    "55+35=",
    "11+78=",
    "56+12=",
    "52-26=",
    "39+12=",
    "6+34=",
    "7+58=",
    "97-53=",
    "50-46=",
    "66-59=",
    "16-2=",
    "11-10=",
    "58+27=",
    "92-29=",
    "21+61=",
    "3+42=",
    "4+1=",
    "36+23=",
    "41+56=",
    "23-20=",
    "84+3=",
    "18+37=",
    "85-49=",
    "52-3=",
    "62-18=",
    "66-46=",
    "11+54=",
    "60+15=",
    "60+34=",
    "9+86=",
    "72-25=",
    "44+17=",
    "33-15=",
    "2+60=",
    "87-2=",
    "17+75=",
    "81-10=",
    "3+94=",
    "36+42=",
    "11+22=",
    "14+38=",
    "93-62=",
    "34-7=",
    "31+26=",
    "30-20=",
    "11+53=",
    "99-20=",
    "44-25=",
    "56-44=",
    "61+9=",
    "56+26=",
    "66-60=",
    "91-42=",
    "70+24=",
    "80+9=",
    "15+11=",
    "48-31=",
    "91-62=",
    "32+18=",
    "71-13=",
    "70+11=",
    "71-66=",
    "60-59=",
    "43+5=",
    "32-12=",
    "95-74=",
    "70+9=",
    "79-39=",
    "75-46=",
    "57-18=",
    "21+7=",
    "87+3="
)

$table = $d.Tables.Item(1)
$rows = $table.Rows.Count
$cols = $table.Columns.Count

$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $table.Cell($r, $c)
        $range = $cell.Range
        # Trim trailing cell-mark / paragraph-mark characters before setting text
        $range.End = $range.End - 1
        $range.Text = $newValues[$idx]
        $idx++
    }
}
